$d = $word.ActiveDocument

$replacements = @(
    @("43×34=", "29×82="),
    @("28×50=", "18×26="),
    @("81×14=", "47×40="),
    @("74×88=", "48×98="),
    @("31×85=", "97×41="),
    @("29×95=", "72×40="),
    @("15×34=", "90×72="),
    @("76×97=", "67×58="),
    @("55×71=", "52×21="),
    @("71×33=", "90×39="),
    @("97×62=", "92×83="),
    @("12×55=", "58×14="),
    @("52×65=", "85×18="),
    @("11×83=", "13×64="),
    @("80×68=", "17×57="),
    @("72×55=", "53×67="),
    @("52×93=", "36×55="),
    @("53×26=", "78×85="),
    @("70×16=", "45×32="),
    @("11×29=", "59×53="),
    @("76×40=", "77×81="),
    @("29×68=", "86×39="),
    @("33×77=", "28×21="),
    @("80×74=", "97×99="),
    @("66×41=", "51×62=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
